# "Örnek 10 - Metin İşlemleri.xlsx" güncellemesi
# Kullanıcı (öğrenci) kendi bilgilerini girip metin işlevlerinin
# sonuçlarını görüyor: A5 hücresine ad-soyad yazılıyor, D5:H5 bu veriden
# türetilen LEN / UPPER / LOWER / MID / DOLLAR sonuçlarını gösteriyor;
# E8:E10'un yanındaki (F sütunu) alanlara da Numara / Ad Soyad / Bölüm
# bilgileri giriliyor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A5:C5 (birleşik hücre) -> ad soyad girişi -------------------------
$ws.Range("A5").Value = "yusufkarataş"

# --- D5:H5 -> metin işlevleri -------------------------------------------
$ws.Range("D5").Formula = "=LEN(A5)"
$ws.Range("E5").Formula = "=UPPER(A5)"
$ws.Range("F5").Formula = "=LOWER(A5)"
$ws.Range("G5").Value = "fkar"

# H5: D5'teki sayıyı (uzunluk) Lira biçiminde göstermek için DOLLAR().
# Hücrenin sayı biçimini de Türk Lirası gösterimine çeviriyoruz.
$ws.Range("H5").NumberFormat = '"₺"#,##0.00'
$ws.Range("H5").Formula = "=DOLLAR(D5)"

# --- Numara / Ad Soyad / Bölüm bilgi kutuları (F8:F10) ------------------
$ws.Range("F8").Value = 20215070055
$ws.Range("F9").Value = "Muhammed Ali Harmancı"
$ws.Range("F10").Value = "Yönetim Bilişim Sistemleri"

# --- E sütunu: büyütülmüş metni (UPPER sonucu) daha iyi göstermek için
#     sütun genişliği biraz artırıldı.
$ws.Columns.Item(5).ColumnWidth = 24.65

# --- Seçili alan A5:C5 olarak bırakılıyor (yeniden adlandırılan alan) --
[void]$ws.Range("A5:C5").Select()
